$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.990.39"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.91"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.59"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4591"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07716"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9797"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.07"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.896.80"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.81%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.673"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.933"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.84"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009460"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.48%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.70"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.954.76"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.320"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = "Toncoin"

$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.093"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "Monero"

$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.14"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "EthereumClassic"

$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.07"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "InternetComputer(DFINITY)"

$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.653"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "BitcoinCash"

$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.55"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "LidoDAOToken"

$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.851"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "Stellar"

$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09284"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "ImmutableX"

$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8646"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "Filecoin"

$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.069"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "ARBITRUM"

$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.246"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.42%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "HuobiToken"

$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.022"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "Hedera"

$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05739"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "TrustWalletToken"

$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.154"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "Frax"

$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "VeChain"

$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02040"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "TheSandbox"

$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5510"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "FraxShare"

$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.409"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "Algorand"

$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1754"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "MXToken"

$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.867"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "Aptos"

$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.334"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "Decentraland"

$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5177"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "EnergySwap"

$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.26"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "Cronos"

$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06837"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "PEPE"

$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000002605"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.045"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.95"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "PaxDollar"

$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.11%  "
$ws.Range("E51").Style = "Normal"
